$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Price (D) / Volume(1h) (E) text values.
# Price cells are forced to text (leading apostrophe) so strings like
# "1.000" / "0.9991" / "29.428.79" are preserved verbatim instead of
# being coerced into numbers by Excel; Style is reset to "Normal" right
# after so the quote-prefix formatting does not leave a stray cell style.
$updates = @(
    @{ Row = 2; D = '29.428.79'; E = '  -0.51%  ' },
    @{ Row = 3; D = '1.851.21'; E = $null },
    @{ Row = 4; D = '0.9991'; E = '  -0.01%  ' },
    @{ Row = 5; D = '241.28'; E = '  -1.11%  ' },
    @{ Row = 6; D = '0.6333'; E = '  -1.66%  ' },
    @{ Row = 7; D = $null; E = '  +0.03%  ' },
    @{ Row = 8; D = '4.656.13'; E = '  +142.99%  ' },
    @{ Row = 9; D = '4.764.27'; E = '  +119.79%  ' },
    @{ Row = 10; D = '0.07572'; E = '  +0.46%  ' },
    @{ Row = 11; D = '0.2960'; E = '  -1.75%  ' },
    @{ Row = 12; D = '24.62'; E = '  +0.93%  ' },
    @{ Row = 13; D = '0.07727'; E = '  +0.77%  ' },
    @{ Row = 14; D = '4.999'; E = '  -1.10%  ' },
    @{ Row = 15; D = '0.6847'; E = $null },
    @{ Row = 16; D = '83.02'; E = '  -1.17%  ' },
    @{ Row = 17; D = '0.000009921'; E = '  +3.36%  ' },
    @{ Row = 18; D = '6.187'; E = '  -1.44%  ' },
    @{ Row = 19; D = '29.459.69'; E = $null },
    @{ Row = 20; D = '231.84'; E = '  -2.42%  ' },
    @{ Row = 21; D = $null; E = '  -1.30%  ' },
    @{ Row = 23; D = '7.610'; E = '  -1.37%  ' },
    @{ Row = 24; D = '1.000'; E = '  -0.01%  ' },
    @{ Row = 25; D = '155.94'; E = $null },
    @{ Row = 26; D = '0.1388'; E = '  -1.98%  ' },
    @{ Row = 27; D = '8.407'; E = '  -1.49%  ' },
    @{ Row = 28; D = $null; E = '  -0.83%  ' },
    @{ Row = 29; D = '4.866.87'; E = '  +134.91%  ' },
    @{ Row = 30; D = '1.470'; E = '  -1.31%  ' },
    @{ Row = 31; D = '0.05762'; E = '  -3.58%  ' },
    @{ Row = 32; D = '1.258'; E = $null },
    @{ Row = 33; D = '4.131'; E = '  -0.48%  ' },
    @{ Row = 34; D = '4.019'; E = '  -1.77%  ' },
    @{ Row = 35; D = $null; E = '  -1.63%  ' },
    @{ Row = 36; D = $null; E = '  -1.46%  ' },
    @{ Row = 37; D = '0.7168'; E = '  -1.05%  ' },
    @{ Row = 38; D = '2.596'; E = '  -0.21%  ' },
    @{ Row = 39; D = '1.254.56'; E = '  +3.54%  ' },
    @{ Row = 40; D = $null; E = '  +0.79%  ' },
    @{ Row = 41; D = '0.01804'; E = '  +1.36%  ' },
    @{ Row = 42; D = '0.9027'; E = '  -1.45%  ' },
    @{ Row = 43; D = '6.130'; E = '  -0.89%  ' },
    @{ Row = 44; D = '0.9998'; E = '  +0.03%  ' },
    @{ Row = 45; D = $null; E = '  -0.30%  ' },
    @{ Row = 46; D = '66.97'; E = '  -0.73%  ' },
    @{ Row = 47; D = '7.149'; E = '  -3.22%  ' },
    @{ Row = 48; D = '9.158'; E = '  -0.46%  ' },
    @{ Row = 49; D = '0.4023'; E = '  -1.16%  ' },
    @{ Row = 50; D = '1.685'; E = '  +1.08%  ' },
    @{ Row = 51; D = '0.1125'; E = '  -0.55%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Leading apostrophe forces Excel to store the Price column as literal
        # text (so "1.000", "0.9991", "29.428.79", etc. keep their exact
        # digits/trailing zeros instead of being parsed as numbers).
        $ws.Range("D$($u.Row)").Value = "`'$($u.D)"
        $ws.Range("D$($u.Row)").Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
